$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fix the header typo ("elease Date" -> "Release Date")
$ws.Range("A1").Value = "Release Date"

# Strip the " (Qn)" quarter suffix from the release-date strings in column A
$dates = @(
    "Dec 02, 2019",
    "Sep 02, 2019",
    "May 31, 2019",
    "Mar 11, 2019",
    "Dec 10, 2018",
    "Sep 10, 2018",
    "Jun 11, 2018",
    "Mar 29, 2018",
    "Dec 11, 2017",
    "Sep 11, 2017",
    "Jun 12, 2017",
    "Mar 31, 2017",
    "Dec 12, 2016",
    "Sep 09, 2016",
    "Jun 10, 2016",
    "Mar 31, 2016",
    "Dec 10, 2015",
    "Sep 10, 2015",
    "Jun 10, 2015",
    "Mar 31, 2015",
    "Dec 10, 2014",
    "Sep 10, 2014"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
}

# Update the active selection to match the saved view state
$ws.Range("C7").Select()
